# Rename the "HMOX1+ prooxidative papillary" cluster label to
# "HMOX+ anti-oxidative papillary" on the subclusters sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fresh only - subclusters")

$ws.Range("D8:D13").Value = "HMOX+ anti-oxidative papillary"

# Leave the selection on the edited range, matching the saved workbook state.
$ws.Range("D8:D13").Select()
